# Task 05 Answers - Q.26: add the missing answer paragraph right after
# the "Q.26 [line 237] What is "hex" and what does it do? (url in your
# notes)" question paragraph.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute(
    "Q.26 [line 237] What is ""hex"" and what does it do? (url in your notes) ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    Write-Host "Could not find the Q.26 question paragraph."
} else {
    $questionPara = $find.Parent.Paragraphs(1)

    # Split a brand new paragraph off right after the question paragraph.
    $questionPara.Range.InsertParagraphAfter()
    $answerPara = $questionPara.Next()

    # Fill it with the answer text.
    $answerPara.Range.Text = "Hex is the memory address of a variable used by pointers."

    # The new paragraph/run inherited the question's bold run formatting
    # and "space after 0" paragraph spacing from the split point - the
    # answer itself is plain body text, so clear that back off.
    $answerPara.Range.Font.Bold = $false
    $answerPara.Range.Font.BoldBi = $false
    $answerPara.Range.ParagraphFormat.SpaceAfter = 8

    Write-Host "Inserted Q.26 answer paragraph."
}
